# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet a new (blank) column is inserted
# before the existing "Late" column so the schedule can carry an extra
# "Variable Instalment" style data point. All cells from the old column N
# onward (Late / heading / Outstanding) shift one column to the right
# (N->O, O->P, P->Q). The newly inserted column N itself is left blank,
# keeping the same width as its neighbour, column M ("Outstanding").
#
# The user then switches from the "Transactions" tab back to the
# "Repayment schedule" tab and leaves the selection on cell K18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N, shifting "Late"/"heading"/"Outstanding"
# (and all their data cells) one column to the right.
$ws.Columns("N").Insert()

# Give the freshly inserted column the same width as column M
# ("Outstanding"), matching the rest of the schedule's formatting.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet (it was "Transactions"
# before), and leave the cursor on K18.
$ws.Activate()
$ws.Range("K18").Select()
